$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells whose new value looks like a plain number (e.g. "606.78") would be
# auto-parsed into a numeric cell by Excel on assignment, losing the literal text
# formatting (trailing zeros, etc). Force Text format on those first so the string
# is preserved verbatim, matching how the source feed stores "Price" as text.
$priceTextCells = @("D5","D6","D14","D19","D20","D21","D22","D23","D24","D25","D28","D30","D31","D38","D41","D43","D44","D45","D46","D47","D48")
foreach ($cellRef in $priceTextCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "63.916.11"
$ws.Range("E2").Value = "  +0.77%  "
$ws.Range("D3").Value = "3.335.62"
$ws.Range("E3").Value = "  +2.63%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "606.78"
$ws.Range("E5").Value = "  +2.16%  "
$ws.Range("D6").Value = "143.37"
$ws.Range("E6").Value = "  +1.22%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "3.334.34"
$ws.Range("E8").Value = "  +2.63%  "
$ws.Range("E9").Value = "  +0.20%  "
$ws.Range("E10").Value = "  +1.91%  "
$ws.Range("E11").Value = "  +4.22%  "
$ws.Range("E12").Value = "  +1.31%  "
$ws.Range("E13").Value = "  +0.75%  "
$ws.Range("D14").Value = "35.28"
$ws.Range("E14").Value = "  +2.68%  "
$ws.Range("D15").Value = "3.886.86"
$ws.Range("E15").Value = "  +2.74%  "
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("D17").Value = "3.338.81"
$ws.Range("E17").Value = "  +2.67%  "
$ws.Range("D18").Value = "64.031.39"
$ws.Range("E18").Value = "  +1.08%  "
$ws.Range("D19").Value = "6.89"
$ws.Range("E19").Value = "  +1.71%  "
$ws.Range("D20").Value = "483.88"
$ws.Range("E20").Value = "  +1.60%  "
$ws.Range("D21").Value = "14.16"
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("D22").Value = "0.741"
$ws.Range("E22").Value = "  +2.47%  "
$ws.Range("D23").Value = "8.00"
$ws.Range("E23").Value = "  +0.75%  "
$ws.Range("D24").Value = "14.07"
$ws.Range("E24").Value = "  +6.54%  "
$ws.Range("D25").Value = "85.11"
$ws.Range("E25").Value = "  +1.21%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("E27").Value = "  +2.03%  "
$ws.Range("D28").Value = "8.33"
$ws.Range("E28").Value = "  +3.34%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("D30").Value = "7.23"
$ws.Range("E30").Value = "  -3.25%  "
$ws.Range("D31").Value = "2.17"
$ws.Range("E31").Value = "  +2.25%  "
$ws.Range("E32").Value = "  +5.04%  "
$ws.Range("E33").Value = "  -1.11%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("E35").Value = "  +1.82%  "
$ws.Range("E36").Value = "  +3.61%  "
$ws.Range("D37").Value = "0.0₃0756"
$ws.Range("E37").Value = "  +5.85%  "
$ws.Range("D38").Value = "52.51"
$ws.Range("E38").Value = "  -0.69%  "
$ws.Range("E39").Value = "  +2.19%  "
$ws.Range("D40").Value = "3.141.65"
$ws.Range("E40").Value = "  +5.48%  "
$ws.Range("D41").Value = "435.54"
$ws.Range("E41").Value = "  +3.65%  "
$ws.Range("E42").Value = "  +7.66%  "
$ws.Range("D43").Value = "2.78"
$ws.Range("E43").Value = "  +0.92%  "
$ws.Range("D44").Value = "8.40"
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("D45").Value = "0.269"
$ws.Range("E45").Value = "  +0.33%  "
$ws.Range("D46").Value = "2.28"
$ws.Range("E46").Value = "  +5.12%  "
$ws.Range("D47").Value = "37.25"
$ws.Range("E47").Value = "  +10.29%  "
$ws.Range("D48").Value = "26.58"
$ws.Range("E48").Value = "  +2.58%  "
$ws.Range("E50").Value = "  -0.09%  "
$ws.Range("E51").Value = "  -0.41%  "
